$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the "Posted At" column keeps storing plain text dates (not Excel date serials)
$ws.Range("F2:F3").NumberFormat = "@"

# Update row 2 with new job match data
$ws.Range("A2").Value = "Sr Data Engineer"
$ws.Range("B2").Value = "nan"
$ws.Range("C2").Value = "Glendale, CA, US USA"
$ws.Range("D2").Value = 18.9
$ws.Range("E2").Value = "AI Engineer, Data Scientist, LangChain, RAG, FAISS, Pinecone, S3, Glue, Kinesis, Databricks"
$ws.Range("F2").Value = "2026-02-27"
$ws.Range("G2").Value = "https://www.indeed.com/viewjob?jk=6f61d4f9e0de29b6"

# Update row 3 with new job match data (location unchanged)
$ws.Range("A3").Value = "Site Reliability Engineer II"
$ws.Range("B3").Value = "nan"
$ws.Range("C3").Value = "New York, NY, US USA"
$ws.Range("D3").Value = 15.6
$ws.Range("E3").Value = "Data Scientist, Copilot, Docker, Kubernetes, CI/CD, Jenkins, GitHub Actions, Terraform, Git, Python"
$ws.Range("F3").Value = "2026-02-27"
$ws.Range("G3").Value = "https://www.indeed.com/viewjob?jk=4420b6c7bc5a59dc"

# Remove old rows 4-7 which are no longer part of the match list
$ws.Range("A4:G7").EntireRow.Delete()

$wb.Save()
